# Table of Contents update:
#  - Remove the three "Ash Grove" violin part rows (old rows 53-55)
#  - Add new entries for an Irish Lament, a set of Slides, and a set of
#    Slip Jigs (new rows 53-59), renumbering the page column accordingly
#  - "Gravel Walk" entry itself is unchanged in content, but the commit
#    message references updating that section's surrounding version/page
#    numbers, which falls out naturally from the row shuffle above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "The Ash Grove" rows (53: Violin 1, 54: Violin 2, 55: violin 3)
$ws.Rows("53:55").Delete()

# Row 53: Irish Lament
$ws.Cells.Item(53, 1).Value2 = "The Foggy Dew"
$ws.Cells.Item(53, 2).Value2 = "Irish Lament"
$ws.Cells.Item(53, 3).Value2 = 26

# Rows 54-56: Slide set (enter the shared "Slide" type before each tune name,
# matching the order tunes were typed in when the sheet was authored)
$ws.Cells.Item(54, 2).Value2 = "Slide"
$ws.Cells.Item(54, 1).Value2 = "Road to Lisdoonvarna"
$ws.Cells.Item(54, 3).Value2 = 27

$ws.Cells.Item(55, 2).Value2 = "Slide"
$ws.Cells.Item(55, 1).Value2 = "Mick Duggan's Slide"
$ws.Cells.Item(55, 3).Value2 = 27

$ws.Cells.Item(56, 2).Value2 = "Slide"
$ws.Cells.Item(56, 1).Value2 = "Denis Murphy's"
$ws.Cells.Item(56, 3).Value2 = 27

# Rows 57-59: Slip Jig set (tune name entered before the shared "Slip Jig " type)
$ws.Cells.Item(57, 1).Value2 = "A Fig for a Kiss"
$ws.Cells.Item(57, 2).Value2 = "Slip Jig "
$ws.Cells.Item(57, 3).Value2 = 28

$ws.Cells.Item(58, 1).Value2 = "The Butterfly"
$ws.Cells.Item(58, 2).Value2 = "Slip Jig "
$ws.Cells.Item(58, 3).Value2 = 28

$ws.Cells.Item(59, 1).Value2 = "Rocky Road to Dublin"
$ws.Cells.Item(59, 2).Value2 = "Slip Jig "
$ws.Cells.Item(59, 3).Value2 = 28

# Update the view/selection to reflect the newly added bottom of the list
$ws.Range("A66").Select()
